$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-24 Wednesday" "2024-01-25 Thursday"

Replace-Text "86÷2=" "70÷3="
Replace-Text "11÷4=" "66÷3="
Replace-Text "86÷5=" "48÷2="
Replace-Text "18÷9=" "96÷7="
Replace-Text "85÷8=" "39÷7="
Replace-Text "80÷3=" "80÷7="
Replace-Text "49÷2=" "86÷7="
Replace-Text "43÷5=" "43÷8="
Replace-Text "13÷7=" "64÷6="
Replace-Text "31÷3=" "32÷3="
Replace-Text "33÷8=" "46÷8="
Replace-Text "98÷6=" "40÷4="
Replace-Text "85÷6=" "64÷4="
Replace-Text "83÷6=" "19÷5="
Replace-Text "59÷3=" "66÷3="
Replace-Text "92÷5=" "38÷4="
Replace-Text "29÷7=" "82÷3="
Replace-Text "32÷9=" "91÷8="
Replace-Text "38÷6=" "52÷5="
Replace-Text "60÷4=" "46÷4="
Replace-Text "66÷6=" "35÷8="
Replace-Text "93÷5=" "45÷3="
Replace-Text "16÷2=" "85÷3="
Replace-Text "33÷2=" "96÷4="
Replace-Text "80÷6=" "57÷6="

Write-Host "Done applying replacements"
